$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Map of (row, col) -> new text for the arithmetic-problem table.
# Row numbers correspond to the Word table's 1-based row index
# (rows with actual content are 1, 5, 10, 15, 20 — the others are
# blank spacer rows left untouched).
$updates = @(
    @{ Row = 1;  Col = 1; New = "283×5=" },
    @{ Row = 1;  Col = 2; New = "819×2=" },
    @{ Row = 1;  Col = 3; New = "907×3=" },
    @{ Row = 1;  Col = 4; New = "471×6=" },
    @{ Row = 1;  Col = 5; New = "360×5=" },

    @{ Row = 5;  Col = 1; New = "365×6=" },
    @{ Row = 5;  Col = 2; New = "360×2=" },
    @{ Row = 5;  Col = 3; New = "362×3=" },
    @{ Row = 5;  Col = 4; New = "397×3=" },
    @{ Row = 5;  Col = 5; New = "280×4=" },

    @{ Row = 10; Col = 1; New = "164×9=" },
    @{ Row = 10; Col = 2; New = "284×5=" },
    @{ Row = 10; Col = 3; New = "782×5=" },
    @{ Row = 10; Col = 4; New = "884×7=" },
    @{ Row = 10; Col = 5; New = "994×8=" },

    @{ Row = 15; Col = 1; New = "835×3=" },
    @{ Row = 15; Col = 2; New = "595×4=" },
    @{ Row = 15; Col = 3; New = "733×4=" },
    @{ Row = 15; Col = 4; New = "779×7=" },
    @{ Row = 15; Col = 5; New = "191×5=" },

    @{ Row = 20; Col = 1; New = "423×6=" },
    @{ Row = 20; Col = 2; New = "258×4=" },
    @{ Row = 20; Col = 3; New = "513×2=" },
    @{ Row = 20; Col = 4; New = "723×9=" },
    @{ Row = 20; Col = 5; New = "746×5=" }
)

foreach ($u in $updates) {
    $cell = $tbl.Rows.Item($u.Row).Cells.Item($u.Col)
    $cell.Range.Text = $u.New
}

Write-Host "Updated $($updates.Count) cells"
